# CloseBook-RollOver data changes:
# - add script/test-data for Accounts and Closebook-Rollover modules
#   (new "Negative.closebook" scenario row) and correct the org state
#   value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# orgstate value correction: "Goa" -> "Bangalore"
$ws.Range("B19").Value = "Bangalore"

# new test-data row added after "newpassword"/"testclosebook@123"
$ws.Range("A51").Value = "Negative.closebook"
$ws.Range("B51").Value = "31-02-2024,02-31-2024,123, ,,string123,@,etc"

# widen column B so the new, longer value is readable
$ws.Columns("B").ColumnWidth = 43.57142857142857

# reflect the author's final cursor position/selection when editing
$app = $ws.Application
$win = $app.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
[void]$ws.Range("B30").Select()
